# Teacher availability constraint added
# Update the timetable cells to reflect the new teacher availability.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 19 (1:00 PM - 2:00 PM slot)
$ws.Range("E19").Value = "DBMS"
$ws.Range("F19").Value = "DBMS"

# Row 22 (2:00 PM - 3:00 PM slot)
$ws.Range("B22").Value = "DBMS"
$ws.Range("D22").Value = "AI"
$ws.Range("E22").Value = "CO"

# Row 25 (3:00 PM - 4:00 PM slot)
$ws.Range("B25").Value = "AI"
$ws.Range("C25").Value = "CO"
